$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-19 Friday" "2025-12-20 Saturday"

Replace-Text "33×37=1221" "78×51=3978"
Replace-Text "54×25=1350" "85×67=5695"
Replace-Text "14×49=686" "36×38=1368"
Replace-Text "35×18=630" "40×42=1680"
Replace-Text "44×98=4312" "14×80=1120"

Replace-Text "74×34=2516" "16×33=528"
Replace-Text "82×61=5002" "20×98=1960"
Replace-Text "13×98=1274" "17×77=1309"
Replace-Text "96×72=6912" "81×72=5832"
Replace-Text "46×86=3956" "54×15=810"

Replace-Text "75×34=2550" "30×57=1710"
Replace-Text "87×20=1740" "96×78=7488"
Replace-Text "14×15=210" "55×78=4290"
Replace-Text "49×27=1323" "97×27=2619"
Replace-Text "28×40=1120" "49×75=3675"

Replace-Text "90×53=4770" "17×75=1275"
Replace-Text "67×11=737" "95×28=2660"
Replace-Text "55×64=3520" "75×13=975"
Replace-Text "63×93=5859" "30×43=1290"
Replace-Text "48×56=2688" "57×32=1824"

Replace-Text "92×24=2208" "50×90=4500"
Replace-Text "60×58=3480" "72×26=1872"
Replace-Text "60×29=1740" "73×37=2701"
Replace-Text "98×36=3528" "46×82=3772"
Replace-Text "79×57=4503" "96×85=8160"
